$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.971104264259338
$ws.Range("B1").Value = 2.576746940612793
$ws.Range("C1").Value = 2.719305753707886
$ws.Range("D1").Value = 3.385390043258667
$ws.Range("E1").Value = 1.036555886268616
